# Applies the diff: row 16 column B value update, and rows 17/18 content swap
# (Id/Taxonsorteringsordning/Rodlistade/TaxonId/Artnamn/Vetenskapligt
# namn/Auktor/Alder-Stadium/Kon/Aktivitet/Metod/Ost/Nord) with new, independent
# Taxonsorteringsordning (column B) values on both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: only Taxonsorteringsordning (B) changes ---
$ws.Range("B16").Value = 81697

# --- Preserve the K:N block (Alder-Stadium/Kon/Aktivitet/Metod) from row 17,
#     which is the only part of the row best handled via copy (it creates/
#     removes blank cells cleanly), then clear it from row 17 since the new
#     row 17 content (old row 18) does not have that block. ---
$ws.Range("K17:N17").Copy($ws.Range("K18:N18"))
$ws.Range("K17:N17").ClearContents()

# --- Row 17 becomes the former row 18 record, with a new Taxonsorteringsordning ---
$ws.Range("A17").Value = 111939897
$ws.Range("B17").Value = 98961
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 222498
$ws.Range("F17").Value = "Blåsippa"
$ws.Range("G17").Value = "Hepatica nobilis"
$ws.Range("H17").Value = "Schreb."
$ws.Range("Q17").Value = 653206
$ws.Range("R17").Value = 6599944

# --- Row 18 becomes the former row 17 record, with a new Taxonsorteringsordning ---
$ws.Range("A18").Value = 111939910
$ws.Range("B18").Value = 56446
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 100049
$ws.Range("F18").Value = "Spillkråka"
$ws.Range("G18").Value = "Dryocopus martius"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("Q18").Value = 653148
$ws.Range("R18").Value = 6600341
